$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save original values for rows 3-6, columns C (runs), D (balls), E (fours)
# (.Value2 is used for reads/writes here - .Value is unreliable for round-tripping
# a COM Variant in this host, it stringifies to the property descriptor instead
# of the underlying data)
$orig3 = @($ws.Range("C3").Value2, $ws.Range("D3").Value2, $ws.Range("E3").Value2)
$orig4 = @($ws.Range("C4").Value2, $ws.Range("D4").Value2, $ws.Range("E4").Value2)
$orig5 = @($ws.Range("C5").Value2, $ws.Range("D5").Value2, $ws.Range("E5").Value2)
$orig6 = @($ws.Range("C6").Value2, $ws.Range("D6").Value2, $ws.Range("E6").Value2)

# New row3 <- old row6, new row4 <- old row5, new row5 <- old row3, new row6 <- old row4
$ws.Range("C3").Value2 = $orig6[0]
$ws.Range("D3").Value2 = $orig6[1]
$ws.Range("E3").Value2 = $orig6[2]

$ws.Range("C4").Value2 = $orig5[0]
$ws.Range("D4").Value2 = $orig5[1]
$ws.Range("E4").Value2 = $orig5[2]

$ws.Range("C5").Value2 = $orig3[0]
$ws.Range("D5").Value2 = $orig3[1]
$ws.Range("E5").Value2 = $orig3[2]

$ws.Range("C6").Value2 = $orig4[0]
$ws.Range("D6").Value2 = $orig4[1]
$ws.Range("E6").Value2 = $orig4[2]
